# Generate Report for Handback
# ------------------------------------------------------------
# This script reproduces, via Excel COM-interop calls, the "handback
# report generation" edit: it fills in the previously-empty "Latest
# Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on the per-locale sheets (zh-cn, de-de), links the target
# file cells back to their source .md file on GitHub (same as column
# A), flips the overview Status text from "Ready for handoff" to
# "Handed back: in sync with en-US", and widens a few columns that now
# hold longer text.

$wb = $excel.ActiveWorkbook

# Blue/underlined "hyperlink" look used throughout this workbook
# (explicit RGB FF6495ED, matching the existing custom "HyperLink"
# cell style - not the generic themed one).
$HyperlinkColor = 15570276   # BGR-packed 0x64,0x95,0xED == RGB(100,149,237)

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $HyperlinkColor
}

# ----------------------------------------------------------------
# 1. Overview sheet: status text + wider E/F columns
# ----------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"

$ov.Columns.Item(5).ColumnWidth = 29.166666666666668
$ov.Columns.Item(6).ColumnWidth = 29.166666666666668

# ----------------------------------------------------------------
# 2. zh-cn sheet
# ----------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"

$zh.Range("J2").Value = "1ab82184-bb26-4abb-ad0b-a77c84db3a10.ef29714b64eeedcd3b307321aad5f188e0535abb.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-01 02:35:57"
$zh.Range("J3").Value = "2e62452b-dda7-47f8-a59c-cb01d6aad302.d25ae111ca7fe14a73d8dc16ddfae683826e05b4.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-01 02:35:57"

# Rebuild the hyperlink list in document order (A2, I2, A3, I3) so the
# newly-added links for column I slot in right after their row's
# column-A link, matching how Excel numbers r:id relationships.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afeac0f233e373aedd419790ffd86ef5d2017ba1/e2e/1ab82184-bb26-4abb-ad0b-a77c84db3a10.md", "", "", "1ab82184-bb26-4abb-ad0b-a77c84db3a10.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afeac0f233e373aedd419790ffd86ef5d2017ba1/e2e/1ab82184-bb26-4abb-ad0b-a77c84db3a10.md", "", "", "1ab82184-bb26-4abb-ad0b-a77c84db3a10.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afeac0f233e373aedd419790ffd86ef5d2017ba1/e2e/2e62452b-dda7-47f8-a59c-cb01d6aad302.md", "", "", "2e62452b-dda7-47f8-a59c-cb01d6aad302.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afeac0f233e373aedd419790ffd86ef5d2017ba1/e2e/2e62452b-dda7-47f8-a59c-cb01d6aad302.md", "", "", "2e62452b-dda7-47f8-a59c-cb01d6aad302.md")

Style-AsHyperlink $zh.Range("A2")
Style-AsHyperlink $zh.Range("I2")
Style-AsHyperlink $zh.Range("A3")
Style-AsHyperlink $zh.Range("I3")

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ----------------------------------------------------------------
# 3. de-de sheet
# ----------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

$de.Range("J2").Value = "1ab82184-bb26-4abb-ad0b-a77c84db3a10.ef29714b64eeedcd3b307321aad5f188e0535abb.de-de.xlf"
$de.Range("K2").Value = "2016-09-01 02:36:12"
$de.Range("J3").Value = "2e62452b-dda7-47f8-a59c-cb01d6aad302.d25ae111ca7fe14a73d8dc16ddfae683826e05b4.de-de.xlf"
$de.Range("K3").Value = "2016-09-01 02:36:12"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afeac0f233e373aedd419790ffd86ef5d2017ba1/e2e/1ab82184-bb26-4abb-ad0b-a77c84db3a10.md", "", "", "1ab82184-bb26-4abb-ad0b-a77c84db3a10.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afeac0f233e373aedd419790ffd86ef5d2017ba1/e2e/1ab82184-bb26-4abb-ad0b-a77c84db3a10.md", "", "", "1ab82184-bb26-4abb-ad0b-a77c84db3a10.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afeac0f233e373aedd419790ffd86ef5d2017ba1/e2e/2e62452b-dda7-47f8-a59c-cb01d6aad302.md", "", "", "2e62452b-dda7-47f8-a59c-cb01d6aad302.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afeac0f233e373aedd419790ffd86ef5d2017ba1/e2e/2e62452b-dda7-47f8-a59c-cb01d6aad302.md", "", "", "2e62452b-dda7-47f8-a59c-cb01d6aad302.md")

Style-AsHyperlink $de.Range("A2")
Style-AsHyperlink $de.Range("I2")
Style-AsHyperlink $de.Range("A3")
Style-AsHyperlink $de.Range("I3")

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664
